$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("C11").Value = "Eric Martel"
$ws.Range("D11").Value = 988639
$ws.Range("E11").Value = 6.795
$ws.Range("G11").Value = 19
$ws.Range("H11").Value = 1629
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = 0.7557
$ws.Range("L11").Value = 1629
$ws.Range("M11").Value = 9
$ws.Range("P11").Value = 11.111111111111
$ws.Range("V11").Value = 1
$ws.Range("W11").Value = 0
$ws.Range("Z11").Value = 2
$ws.Range("AA11").Value = 0.9062966
$ws.Range("AB11").Value = 1210
$ws.Range("AC11").Value = 4
$ws.Range("AD11").Value = 11
$ws.Range("AE11").Value = 804
$ws.Range("AF11").Value = 86.918918918919
$ws.Range("AG11").Value = 925
$ws.Range("AH11").Value = 514
$ws.Range("AI11").Value = 290
$ws.Range("AJ11").Value = 92
$ws.Range("AK11").Value = 18
$ws.Range("AL11").Value = 33.333333333333
$ws.Range("AM11").Value = 2
$ws.Range("AN11").Value = 40
$ws.Range("AO11").Value = 23
$ws.Range("AP11").Value = 22
$ws.Range("AQ11").Value = 3
$ws.Range("AR11").Value = 63
$ws.Range("AS11").Value = 15
$ws.Range("AT11").Value = 95
$ws.Range("AU11").Value = 5
$ws.Range("AY11").Value = 2
$ws.Range("AZ11").Value = 66.666666666667
$ws.Range("BA11").Value = 117
$ws.Range("BB11").Value = 53.181818181818
$ws.Range("BC11").Value = 51
$ws.Range("BD11").Value = 48.571428571429
$ws.Range("BE11").Value = 66
$ws.Range("BF11").Value = 57.391304347826
$ws.Range("BG11").Value = 146
$ws.Range("BH11").Value = 27
$ws.Range("BI11").Value = 31
$ws.Range("BJ11").Value = 0
$ws.Range("BK11").Value = 6
$ws.Range("BN11").Value = 3
$ws.Range("BO11").Value = 29
$ws.Range("BQ11").Value = 135.9
$ws.Range("BS11").Value = 3
$ws.Range("BT11").Value = 121
$ws.Range("BY11").Value = 0
$ws.Range("BZ11").Value = 6
$ws.Range("CB11").Value = 12
$ws.Range("CC11").Value = 27
$ws.Range("CD11").Value = 13
$ws.Range("CF11").Value = 0
$ws.Range("CJ11").Value = 26
$ws.Range("CR11").Value = 11
$ws.Range("CS11").Value = 3
$ws.Range("CT11").Value = 5
$ws.Range("CU11").Value = 103
$ws.Range("CV11").Value = 49
$ws.Range("CZ11").Value = 54
$ws.Range("DA11").Value = 12
$ws.Range("DB11").Value = 54.545454545455
$ws.Range("DF11").Value = 564
$ws.Range("DG11").Value = 361
$ws.Range("DH11").Value = 1
$ws.Range("DI11").Value = 15
$ws.Range("DJ11").Value = 2168780

# Row 12
$ws.Range("C12").Value = "Ísak Bergmann Jóhannesson"
$ws.Range("D12").Value = 1112327
$ws.Range("E12").Value = 6.785
$ws.Range("G12").Value = 14
$ws.Range("H12").Value = 1186
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0.6483
$ws.Range("L12").Value = 1186
$ws.Range("M12").Value = 5
$ws.Range("P12").Value = 20
$ws.Range("V12").Value = 0
$ws.Range("W12").Value = 1
$ws.Range("Z12").Value = 1
$ws.Range("AA12").Value = 2.59641267
$ws.Range("AB12").Value = 722
$ws.Range("AC12").Value = 1
$ws.Range("AD12").Value = 21
$ws.Range("AE12").Value = 473
$ws.Range("AF12").Value = 84.014209591474
$ws.Range("AG12").Value = 563
$ws.Range("AH12").Value = 206
$ws.Range("AI12").Value = 267
$ws.Range("AJ12").Value = 139
$ws.Range("AK12").Value = 34
$ws.Range("AL12").Value = 45.945945945946
$ws.Range("AM12").Value = 17
$ws.Range("AN12").Value = 33.333333333333
$ws.Range("AO12").Value = 3
$ws.Range("AP12").Value = 18
$ws.Range("AQ12").Value = 1
$ws.Range("AR12").Value = 43
$ws.Range("AS12").Value = 14
$ws.Range("AT12").Value = 12
$ws.Range("AU12").Value = 1
$ws.Range("AY12").Value = 3
$ws.Range("AZ12").Value = 75
$ws.Range("BA12").Value = 31
$ws.Range("BB12").Value = 43.055555555556
$ws.Range("BC12").Value = 25
$ws.Range("BD12").Value = 45.454545454545
$ws.Range("BE12").Value = 6
$ws.Range("BF12").Value = 35.294117647059
$ws.Range("BG12").Value = 149
$ws.Range("BH12").Value = 8
$ws.Range("BI12").Value = 4
$ws.Range("BJ12").Value = 1
$ws.Range("BK12").Value = 3
$ws.Range("BN12").Value = 0
$ws.Range("BO12").Value = 23
$ws.Range("BQ12").Value = 135.7
$ws.Range("BS12").Value = 2
$ws.Range("BT12").Value = 90
$ws.Range("BY12").Value = 3
$ws.Range("BZ12").Value = 2
$ws.Range("CB12").Value = 7
$ws.Range("CC12").Value = 44
$ws.Range("CD12").Value = 26
$ws.Range("CF12").Value = 1
$ws.Range("CJ12").Value = 20
$ws.Range("CR12").Value = 20
$ws.Range("CS12").Value = 4
$ws.Range("CT12").Value = 51
$ws.Range("CU12").Value = 41
$ws.Range("CV12").Value = 11
$ws.Range("CZ12").Value = 74
$ws.Range("DA12").Value = 10
$ws.Range("DB12").Value = 55.555555555556
$ws.Range("DF12").Value = 222
$ws.Range("DG12").Value = 341
$ws.Range("DH12").Value = 0
$ws.Range("DI12").Value = 5
$ws.Range("DJ12").Value = 2168750
